# Applies the "Updated symbol list" edit described in the task diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "242.66"
$ws.Range("D4").Value = "5.218"
$ws.Range("D5").Value = "0.05604"
$ws.Range("D6").Value = "3.367"
$ws.Range("D7").Value = "6.377"
$ws.Range("D8").Value = "0.8052"
$ws.Range("D9").Value = "0.9541"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.01114"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1437"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07292"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03151"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03095"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09282"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.573"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001663"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
# 0.04700 would lose its trailing zero if stored as a number; force text
# (matches the source data's inline-string type) with a leading apostrophe,
# same as typing '0.04700 directly into the cell.
$ws.Range("D18").Value = "'0.04700"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "0.006347"
$ws.Range("D20").Value = "0.004986"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D25").Value = "2.094"
$ws.Range("D26").Value = "0.3268"
$ws.Range("D40").Value = "0.03914"
$ws.Range("D41").Value = "0.006896"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003400"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "0.1033"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "0.007514"
$ws.Range("D45").Value = "0.00005942"
$ws.Range("D47").Value = "0.0005497"
$ws.Range("D48").Value = "0.6825"
$ws.Range("D49").Value = "'0.07660"
$ws.Range("E49").Value = "48BOLOBOLO"
